# Fixed quantities in BOM after component change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: 12k Resistor (R1,R2,R26,R28,R34,R36,R33) Qty 6 -> 7
$ws.Range("E13").Value = 7

# Row 14: 20k Resistor (R3,R5,R6,R8,R10,R13,R23,R25,R29) Qty 10 -> 9
$ws.Range("E14").Value = 9

# Active cell / selection moved to B13
$ws.Range("B13").Select()
